$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 15 so the old rows 15/16/17 (credits footer block)
# shift down to 16/17/18, opening up rows 13/14 for the two new music
# credit entries.
$ws.Rows("15:15").Insert()

# Fill in the two new rows with the new music files' credit info.
# Column A first (both rows), then column B (both rows), then column C
# (both rows) -- matches the shared-string insertion order seen in the
# target file (new unique strings appended as: A13, A14, B13, B14).
$ws.Range("A13").Value = "Komiku_-_02_-_Boss_4__Cobblestone_in_their_face.mp3"
$ws.Range("A14").Value = "Loyalty_Freak_Music_-_04_-_Cant_Stop_My_Feet_.mp3"
$ws.Range("B13").Value = "https://www.chosic.com/download-audio/25453/"
$ws.Range("B14").Value = "https://www.chosic.com/download-audio/25495/"
$ws.Range("C13").Value = "You are free to use this music in your projects with no required crediting. However, linking back is greatly appreciated. You can use the following text"
$ws.Range("C14").Value = "You are free to use this music in your projects with no required crediting. However, linking back is greatly appreciated. You can use the following text"

# Match the licensing-note formatting used by the other music rows
# (italic Helvetica note style) by copying the format from C11.
$ws.Range("C11").Copy()
$ws.Range("C13:C14").PasteSpecial(-4122)

# The worksheet's Hyperlinks collection doesn't renumber its stored
# ranges when rows are inserted above them, so rebuild it explicitly:
# drop everything and re-add each hyperlink at its (possibly shifted)
# location, in original order, so relationship ids line up again.
$ws.Range("A1").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B6"), "https://opengameart.org/content/spikes-0")
$ws.Hyperlinks.Add($ws.Range("B4"), "https://www.fontspace.com/a-area-kilometer-50-font-f53888")
$ws.Hyperlinks.Add($ws.Range("B2"), "https://freesound.org/people/Whiprealgood/sounds/87535/")
$ws.Hyperlinks.Add($ws.Range("B3"), "https://freesound.org/people/suntemple/sounds/253172/")
$ws.Hyperlinks.Add($ws.Range("B5"), "https://opengameart.org/content/simple-explosion-bleeds-game-art")
$ws.Hyperlinks.Add($ws.Range("B7"), "https://opengameart.org/content/various-inventory-24-pixel-icon-set")
$ws.Hyperlinks.Add($ws.Range("B17"), "https://elthen.itch.io/2d-pixel-art-vegetable-monsters-sprite-pack")
$ws.Hyperlinks.Add($ws.Range("B18"), "https://free-game-assets.itch.io/night-city-street-2d-background-tiles")
$ws.Hyperlinks.Add($ws.Range("B8"), "https://opengameart.org/content/energy-icon")

# Re-adding hyperlinks re-applies hyperlink formatting with a freshly
# interned style; restore the original shared "Hyperlink" cell style on
# every affected cell so formatting matches exactly.
$ws.Range("B2").Style = "Hyperlink"
$ws.Range("B3").Style = "Hyperlink"
$ws.Range("B4").Style = "Hyperlink"
$ws.Range("B5").Style = "Hyperlink"
$ws.Range("B6").Style = "Hyperlink"
$ws.Range("B7").Style = "Hyperlink"
$ws.Range("B8").Style = "Hyperlink"
$ws.Range("B17").Style = "Hyperlink"
$ws.Range("B18").Style = "Hyperlink"

# Match the final selection left behind in the saved file.
$null = $ws.Range("A15").Select()
